# Insert a new data row before the current row 255 (Vega Central Mapocho de
# Santiago - Mango sheet). This pushes the existing rows 255-294 down to
# 256-295, preserving all of their values/formatting, and the newly freed
# row 255 gets the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255, shifting rows 255:294 down to 256:295.
$ws.Rows.Item(255).Insert()

# Populate the new row 255 with the new record's values.
$ws.Cells.Item(255, 1).Value = 9
$ws.Cells.Item(255, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(255, 3).Value = "Metropolitana"
$ws.Cells.Item(255, 4).Value = 44522
$ws.Cells.Item(255, 5).Value = 13
$ws.Cells.Item(255, 6).Value = "Fruta"
$ws.Cells.Item(255, 7).Value = 100108
$ws.Cells.Item(255, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(255, 9).Value = 100108002
$ws.Cells.Item(255, 10).Value = "Mango"
$ws.Cells.Item(255, 11).Value = "Sin especificar"
$ws.Cells.Item(255, 12).Value = "Primera"
$ws.Cells.Item(255, 13).Value = 300
$ws.Cells.Item(255, 14).Value = 6500
$ws.Cells.Item(255, 15).Value = 6500
$ws.Cells.Item(255, 16).Value = 6500
$ws.Cells.Item(255, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(255, 18).Value = "Perú"
$ws.Cells.Item(255, 19).Value = 1625
$ws.Cells.Item(255, 20).Value = 4
